$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 10 ("Objetivos:") — B/C held the wrong text (a teacher name that had
#    been pasted into the wrong row). Replace it with the real objectives
#    text. Style ("2"/"3") and row height are unchanged, so only the value is
#    touched.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "Apresentar aos estudantes de engenharia os conceitos básicos de Ciência dos Materiais."
$ws.Range("C10").Value = "Apresentar aos estudantes de engenharia os conceitos básicos de Ciência dos Materiais."

# ---------------------------------------------------------------------------
# 2) Insert three new rows right after row 12 ("Docentes responsáveis:") to
#    hold the three professors, one per row, in columns B/C only (no label in
#    column A). This pushes the former rows 13-21 down to 16-24.
# ---------------------------------------------------------------------------
$ws.Range("A13:A15").Insert()
$ws.Range("A13:A15").Clear()

$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("B13").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C13").Value = "6495737 - Durval Rodrigues Junior"

$ws.Range("B14").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C14").Value = "5983729 - Fernando Vernilli Junior"

$ws.Range("B15").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C15").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

# ---------------------------------------------------------------------------
# 3) Row 16 ("Programa resumido:") — fill in the short-syllabus text that was
#    previously missing.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "Estrutura e ligação atômica. 2  Estruturas dos materiais. 3  Imperfeições em sólidos. 4  Diagrama de fases. 5  Propriedades mecânicas"
$ws.Range("C16").Value = "Estrutura e ligação atômica. 2  Estruturas dos materiais. 3  Imperfeições em sólidos. 4  Diagrama de fases. 5  Propriedades mecânicas"

# ---------------------------------------------------------------------------
# 4) Row 18 ("Programa:") — fill in the full syllabus text (previously held a
#    misplaced teacher name).
# ---------------------------------------------------------------------------
$programa = @'
1. Estrutura e ligação atômica: estrutura dos átomos; ligações covalente, iônica, metálica e forças de van der Waals.
2. Estruturas dos materiais: sólidos cristalinos; direções e planos cristalográficos; células unitárias; redes de Bravais; fator de empacotamento; métodos para determinação das estruturas cristalinas; estruturas metálicas, iônicas e moleculares. Estrutura de cerâmicas. Estrutura de polímeros. Sólidos amorfos: vidros e polímeros. Aspectos básicos de materiais compósitos. Exemplos de materiais de engenharia.
3. Imperfeições em sólidos: tipos e formação de defeitos; lacunas; soluções sólidas (intersticial e substitucional); estruturas ordenadas; compostos intermetálicos; discordâncias; movimento de discordâncias; defeitos planares (interfaces). Exemplos práticos.
4. Diagrama de fases: definição de fase; regra de Gibbs; curva de resfriamento; diagramas de equilíbrio de sistemas binários; equilíbrio de formação e decomposição de fases. Exemplos de diagramas de fases relacionados com a microestrutura dos materiais.
5. Conceitos básicos sobre as propriedades mecânicas dos materiais: conceitos de tensão e deformação; propriedades elásticas; deformação plástica; plasticidade e fluxo; materiais não newtonianos; relaxação e fluência; fadiga. Exemplos e casos práticos.
'@
$ws.Range("B18").Value = $programa
$ws.Range("C18").Value = $programa

# ---------------------------------------------------------------------------
# 5) Row 21 ("Método:") — shift the evaluation text up from the old row 19
#    ("Critério:") position where it had been misplaced.
# ---------------------------------------------------------------------------
$ws.Range("B21").Value = "Serão aplicadas duas provas escritas com notas P1 e P2."
$ws.Range("C21").Value = "Serão aplicadas duas provas escritas com notas P1 e P2."

# ---------------------------------------------------------------------------
# 6) Row 22 ("Critério:") — formula text shifted up from the old row 20
#    ("Norma de recuperação:") position.
# ---------------------------------------------------------------------------
$ws.Range("B22").Value = "A nota final NF será calculada pela fórmula: NF=(P1 + P2)/2."
$ws.Range("C22").Value = "A nota final NF será calculada pela fórmula: NF=(P1 + P2)/2."

# ---------------------------------------------------------------------------
# 7) Row 23 ("Norma de recuperação:") — text shifted up from the old row 21
#    ("Bibliografia:") position.
# ---------------------------------------------------------------------------
$ws.Range("B23").Value = "Será aplicada uma prova escrita NR que comporá com a nota final NF a média final após recuperação MF=(NF+NF)/2."
$ws.Range("C23").Value = "Será aplicada uma prova escrita NR que comporá com a nota final NF a média final após recuperação MF=(NF+NF)/2."

# ---------------------------------------------------------------------------
# 8) Row 24 ("Bibliografia:") — previously empty, now filled with the full
#    reading list.
# ---------------------------------------------------------------------------
$bibliografia = @'
1) Askeland, D. R.; Phule, P. P. Ciência e engenharia dos materiais. São Paulo: CENGAGE, 2008.
2) Callister Jr., W. D. Fundamentos da ciência e engenharia de materiais. Rio de Janeiro: LTC Editora, 2006.
3) Callister Jr., W. D. Ciência e engenharia de materiais. Rio de Janeiro: LTC Editora, 2008.
4) Van Vlack, L. H. Princípios de ciência e tecnologia dos materiais. Rio de Janeiro: Editora Campus, 1984.
5) Shackelford, J. E. Ciência dos materiais. São Paulo: Prentice Hall, 2008. 
6) Jastrzebski, Z. D. The nature and properties of engineering materials. Nova Iorque: John Wiley, 1987.
7) Padilha, A. F. Materiais de engenharia: microestrutura e propriedades. São Paulo: Hemus Editora, 1997.
8) Ashby, M. F.; Jones, D. R. H. Engenharia de materiais, 2 vol. Rio de Janeiro: Elsevier Editora, 2007.
'@
$ws.Range("B24").Value = $bibliografia
$ws.Range("C24").Value = $bibliografia

Write-Host "edit complete"
